$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.003837439598571479
$ws.Range("E2").Value = 0.3718167378372521
$ws.Range("G2").Value = 0.2494892361375047
$ws.Range("I2").Value = 0.3687475729068005
$ws.Range("L2").Value = 0.597153
$ws.Range("M2").Value = 0.0822565
$ws.Range("N2").Value = 12.82009457445576
$ws.Range("O2").Value = 3.538068880805355

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.00926038528939491
$ws.Range("B2").Value = 0.04893136040142869
$ws.Range("E2").Value = 0.2216906659485062
$ws.Range("I2").Value = 0.4247747675224997
$ws.Range("L2").Value = 0.1116199591040388
$ws.Range("M2").Value = 0.04737166666666669
$ws.Range("N2").Value = 5.019474122500005
$ws.Range("O2").Value = 2.349265400575166

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.08320518682224694
$ws.Range("B2").Value = 0.02875797743582353
$ws.Range("E2").Value = 0.1707961611009805
$ws.Range("I2").Value = 0.4661536766800574
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04855891666666667
$ws.Range("N2").Value = 8.420331206844089
$ws.Range("O2").Value = 4.913850329222639
